$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 420, pushing the existing rows 420-441 down to 421-442
$ws.Rows(420).Insert()

# Populate the newly inserted row 420 with the new weekly record
$ws.Cells.Item(420, 1).Value2 = 4
$ws.Cells.Item(420, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(420, 3).Value2 = "Los Lagos"
$ws.Cells.Item(420, 4).Value2 = 44706
$ws.Cells.Item(420, 5).Value2 = 10
$ws.Cells.Item(420, 6).Value2 = 100112006
$ws.Cells.Item(420, 7).Value2 = "Repollo"
$ws.Cells.Item(420, 8).Value2 = "Crespo record"
$ws.Cells.Item(420, 9).Value2 = "Primera"
$ws.Cells.Item(420, 10).Value2 = 250
$ws.Cells.Item(420, 11).Value2 = 1800
$ws.Cells.Item(420, 12).Value2 = 1800
$ws.Cells.Item(420, 13).Value2 = 1800
$ws.Cells.Item(420, 14).Value2 = "`$/unidad"
$ws.Cells.Item(420, 15).Value2 = "Región del Maule"
$ws.Cells.Item(420, 16).Value2 = 1800
$ws.Cells.Item(420, 17).Value2 = 1
$ws.Cells.Item(420, 18).Value2 = "Hortaliza"

# Ensure the date cell keeps the date number format used throughout column D
$ws.Cells.Item(420, 4).NumberFormat = $ws.Cells.Item(421, 4).NumberFormat
